# Update the statistics summary table (Table 1) with the re-computed
# values (new case/control split, means/medians, CIs and p-values)
# produced after adding the alluvium & spaghetti plots.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(2,4).Range.Text = "    54"
$t.Cell(2,5).Range.Text = "     46"
$t.Cell(3,3).Range.Text = "0.0931 (0.9880)"
$t.Cell(3,4).Range.Text = "0.2383 (0.9085)"
$t.Cell(3,5).Range.Text = "-0.0774 (1.0584)"
$t.Cell(3,6).Range.Text = " 0.11176"
$t.Cell(3,8).Range.Text = " 0.31999"
$t.Cell(4,3).Range.Text = "0.1195 [-0.5691, 0.7773]"
$t.Cell(4,4).Range.Text = "0.3108 [-0.4239, 0.8035]"
$t.Cell(4,5).Range.Text = "-0.1726 [-0.9317, 0.5454]"
$t.Cell(4,6).Range.Text = " 0.07437"
$t.Cell(4,8).Range.Text = " 0.31999"
$t.Cell(5,3).Range.Text = "0.1195 [-2.0104, 3.0897]"
$t.Cell(5,4).Range.Text = "0.3108 [-2.0104, 1.9425]"
$t.Cell(5,5).Range.Text = "-0.1726 [-1.7418, 3.0897]"
$t.Cell(5,6).Range.Text = " 0.07437"
$t.Cell(5,8).Range.Text = " 0.31999"
$t.Cell(6,3).Range.Text = " 26 (26.000) "
$t.Cell(6,4).Range.Text = "12 (22.222) "
$t.Cell(6,5).Range.Text = "14 (30.435) "
$t.Cell(6,6).Range.Text = " 0.51747"
$t.Cell(6,8).Range.Text = " 0.23175"
$t.Cell(7,3).Range.Text = " 36 (36.000) "
$t.Cell(7,4).Range.Text = "19 (35.185) "
$t.Cell(7,5).Range.Text = "17 (36.957) "
$t.Cell(8,3).Range.Text = " 38 (38.000) "
$t.Cell(8,4).Range.Text = "23 (42.593) "
$t.Cell(8,5).Range.Text = "15 (32.609) "
$t.Cell(9,3).Range.Text = "-0.6643 (4.5488)"
$t.Cell(9,4).Range.Text = "-0.7994 (3.9185)"
$t.Cell(9,5).Range.Text = "-0.5058 (5.2334)"
$t.Cell(9,6).Range.Text = " 0.74940"
$t.Cell(9,8).Range.Text = " 0.06353"
$t.Cell(10,3).Range.Text = "-0.7000 [-3.2495, 1.7772]"
$t.Cell(10,4).Range.Text = "-0.7310 [-3.3088, 1.5202]"
$t.Cell(10,5).Range.Text = "-0.4402 [-3.0656, 2.1625]"
$t.Cell(10,6).Range.Text = " 0.53365"
$t.Cell(10,8).Range.Text = " 0.06353"
$t.Cell(11,3).Range.Text = "-0.7000 [-16.2705, 11.4307]"
$t.Cell(11,4).Range.Text = "-0.7310 [-7.8807, 10.7229]"
$t.Cell(11,5).Range.Text = "-0.4402 [-16.2705, 11.4307]"
$t.Cell(11,6).Range.Text = " 0.53365"
$t.Cell(11,8).Range.Text = " 0.06353"
$t.Cell(12,3).Range.Text = " 60 (60.000) "
$t.Cell(12,4).Range.Text = "36 (66.667) "
$t.Cell(12,5).Range.Text = "24 (52.174) "
$t.Cell(12,6).Range.Text = " 0.20421"
$t.Cell(12,8).Range.Text = " 0.29841"
$t.Cell(13,3).Range.Text = " 40 (40.000) "
$t.Cell(13,4).Range.Text = "18 (33.333) "
$t.Cell(13,5).Range.Text = "22 (47.826) "
